# BOM.xlsx update: "Finished Schematics, Libraries, and BOM"
# Adds 8 new BOM line items (rows 23-30) with part number, description,
# quantity and a Digikey hyperlink in column F, widens column A slightly,
# and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM rows: Part Number (A), Part Description (B), Quantity (C), Link (F)
$rows = @(
    @{ Row = 23; A = "TST-105-01-F-D";      B = "STMLink Connector"; C = 1; Url = "https://www.digikey.com/en/products/detail/samtec-inc/TST-105-01-F-D/9497108" },
    @{ Row = 24; A = "SSW-106-02-TM-S-RA";  B = "UART Connector";    C = 1; Url = "https://www.digikey.com/en/products/detail/samtec-inc/SSW-106-02-TM-S-RA/7891818" },
    @{ Row = 25; A = "TL6330AF200Q";        B = "RST Button";        C = 1; Url = "https://www.digikey.com/en/products/detail/e-switch/TL6330AF200Q/8032037" },
    @{ Row = 26; A = "BLM21PG600SN1D";      B = "Ferrite Bead";      C = 1; Url = "https://www.digikey.com/en/products/detail/murata-electronics/BLM21PG600SN1D/584263" },
    @{ Row = 27; A = "5988191107F";         B = "Blue LED";          C = 1; Url = "https://www.digikey.com/en/products/detail/dialight/5988191107F/1291280" },
    @{ Row = 28; A = "LTW-170TK";           B = "White LED";         C = 1; Url = "https://www.digikey.com/en/products/detail/liteon/LTW-170TK/758704" },
    @{ Row = 29; A = "5988110107F";         B = "Red LED";           C = 1; Url = "https://www.digikey.com/en/products/detail/dialight/5988110107F/1291272" },
    @{ Row = 30; A = "5988170107F";         B = "Green LED";         C = 1; Url = "https://www.digikey.com/en/products/detail/dialight/5988170107F/1291278" }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C

    # Put the URL text in column F, then turn it into a hyperlink (keeps the
    # visible text equal to the target URL, same as the existing rows) and
    # apply the workbook's "Hyperlink" cell style so it matches the rest of
    # column F.
    $ws.Cells.Item($r, 6).Value = $item.Url
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $item.Url) | Out-Null
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

# Column A got a little wider to fit the new part numbers/descriptions.
$ws.Columns.Item(1).ColumnWidth = 22.6

# Selection moved to B33 as the last user action before saving.
$ws.Range("B33").Select() | Out-Null
